$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 342 ("Femacal de La Calera" /
# Apio records). This pushes the old rows 342-347 down to 344-349 and
# leaves two blank rows (342-343) to be filled with the new weekly data.
$ws.Rows("342:343").Insert()

# New row 342: week of 2022-03-08 (serial 44628), "Provincia de Santiago"
$ws.Cells.Item(342, 1).Value = 3
$ws.Cells.Item(342, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(342, 3).Value = "Coquimbo"
$ws.Cells.Item(342, 4).Value = 44628
$ws.Cells.Item(342, 5).Value = 5
$ws.Cells.Item(342, 6).Value = 100112017
$ws.Cells.Item(342, 7).Value = "Apio"
$ws.Cells.Item(342, 8).Value = "Americana (o)"
$ws.Cells.Item(342, 9).Value = "Primera"
$ws.Cells.Item(342, 10).Value = 210
$ws.Cells.Item(342, 11).Value = 9000
$ws.Cells.Item(342, 12).Value = 9500
$ws.Cells.Item(342, 13).Value = 9238
$ws.Cells.Item(342, 14).Value = "$/docena de matas"
$ws.Cells.Item(342, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(342, 16).Value = 1540
$ws.Cells.Item(342, 17).Value = 6
$ws.Cells.Item(342, 18).Value = "Hortaliza"

# New row 343: same week, "Segunda" quality, "Provincia de Santiago"
$ws.Cells.Item(343, 1).Value = 3
$ws.Cells.Item(343, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(343, 3).Value = "Coquimbo"
$ws.Cells.Item(343, 4).Value = 44628
$ws.Cells.Item(343, 5).Value = 5
$ws.Cells.Item(343, 6).Value = 100112017
$ws.Cells.Item(343, 7).Value = "Apio"
$ws.Cells.Item(343, 8).Value = "Americana (o)"
$ws.Cells.Item(343, 9).Value = "Segunda"
$ws.Cells.Item(343, 10).Value = 80
$ws.Cells.Item(343, 11).Value = 7500
$ws.Cells.Item(343, 12).Value = 7500
$ws.Cells.Item(343, 13).Value = 7500
$ws.Cells.Item(343, 14).Value = "$/docena de matas"
$ws.Cells.Item(343, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(343, 16).Value = 1250
$ws.Cells.Item(343, 17).Value = 6
$ws.Cells.Item(343, 18).Value = "Hortaliza"
